# testdangnhap.xlsx - "Add files via upload" edit
#
# Summary of the change being applied:
#  - A2 ("hung") becomes the text value "0387335906" (format as Text so the
#    leading zero survives instead of Excel coercing it to a number).
#  - B2 (numeric 1) becomes the text value "beodeptrai1".
#  - Rows 3 and 4 (which only duplicated the old A2/B2 pattern) are removed.
#  - Column B is given an explicit width to fit the new text.
#  - The active selection in the sheet view moves to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A2 -----------------------------------------------------------
# Force a text number format first so the leading zero in "0387335906" is
# preserved instead of being interpreted/stored as a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "0387335906"

# --- Update B2 ------------------------------------------------------------
$ws.Range("B2").Value = "beodeptrai1"

# --- Remove the now unused rows 3 and 4 -----------------------------------
$ws.Range("A3:B4").EntireRow.Delete()

# --- Give column B an explicit custom width for the new text -------------
$ws.Columns.Item(2).ColumnWidth = 9.7

# --- Move the active selection, matching where the author clicked next ---
$ws.Range("C6").Select()
